$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "DropPack" -> "DropPackList" (moves to column Y),
#     and introduce two new header columns: "MoveType" (W) and "AtkDis" (X).
# Order of assignment matters for shared-string allocation order, so write
# Y1, then W1, then X1 (matches target shared-string table ordering).
$ws.Range("Y1").Value = "DropPackList"
$ws.Range("W1").Value = "MoveType"

$atk = $ws.Range("X1")
$atk.Value = "AtkDis"
$atk.Font.Name = "宋体"
$atk.Font.Size = 11
$atk.Font.Family = 3
$atk.Characters(2, 5).Font.Name = "宋体"
$atk.Characters(2, 5).Font.Size = 11

# --- Data rows: new MoveType (W) and AtkDis (X) values per NPC row.
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 20

$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 20

$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 20

$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 20

$ws.Range("W6").Value = 2
$ws.Range("X6").Value = 20

# --- Column widths for the affected/new columns.
$ws.Columns.Item(22).ColumnWidth = 24.25   # V -> stored width 25
$ws.Columns.Item(23).ColumnWidth = 24.25   # W -> stored width 25
$ws.Columns.Item(24).ColumnWidth = 24.25   # X -> stored width 25
$ws.Columns.Item(25).ColumnWidth = 13.08   # Y -> stored width ~13.875

# --- View state: scroll so column K is left-most, select X10.
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("X10").Select()
